$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by 3 rows: insert blank rows 17-19 below the existing
# data (row 16 is currently the last data row) so they pick up the same
# bordered/bold/centered formatting used by column A throughout the table.
$ws.Range("A16:M16").Copy() | Out-Null
$ws.Range("A17:M19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Averaged intensities were recomputed after adding the spiral sampling
# schemes; every scheme row from "Gaussian-Quadrature" onward now carries
# new numbers and the three spiral schemes are brand new rows. Rewrite
# rows 10-19 in full (HKL index, scheme name, and the 11 averaged values)
# to their final values.

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9671164926603516
$ws.Range("D10").Value = 1.1173012170782
$ws.Range("E10").Value = 0.9646396893408864
$ws.Range("F10").Value = 0.9671164926603516
$ws.Range("G10").Value = 1.06341495830896
$ws.Range("H10").Value = 0.9131063717481164
$ws.Range("I10").Value = 0.9658307839933642
$ws.Range("J10").Value = 1.1173012170782
$ws.Range("K10").Value = 1.040970453209543
$ws.Range("L10").Value = 1.004043472934947
$ws.Range("M10").Value = 0.9985682521883131

# Row 11: Spiral-90deg-10rot-5space (new scheme)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9979763275671436
$ws.Range("D11").Value = 0.9063137559222324
$ws.Range("E11").Value = 1.026894610473181
$ws.Range("F11").Value = 0.9979763275671436
$ws.Range("G11").Value = 0.9373924796992724
$ws.Range("H11").Value = 1.089775418852769
$ws.Range("I11").Value = 1.019581208165181
$ws.Range("J11").Value = 0.9063137559222324
$ws.Range("K11").Value = 0.9666041831977066
$ws.Range("L11").Value = 0.9822902553824251
$ws.Range("M11").Value = 0.9963223001132966

# Row 12: Spiral-90deg-15rot-5space (new scheme)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9977361209357697
$ws.Range("D12").Value = 0.9071464250505598
$ws.Range("E12").Value = 1.026721218279455
$ws.Range("F12").Value = 0.9977361209357697
$ws.Range("G12").Value = 0.9378208138508721
$ws.Range("H12").Value = 1.089351098758683
$ws.Range("I12").Value = 1.019401977487226
$ws.Range("J12").Value = 0.9071464250505598
$ws.Range("K12").Value = 0.9669338216650076
$ws.Range("L12").Value = 0.9823349713003886
$ws.Range("M12").Value = 0.9963629423937609

# Row 13: Spiral-90deg-10rot-3space (new scheme)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9979078733253837
$ws.Range("D13").Value = 0.9066257400676264
$ws.Range("E13").Value = 1.026802818205841
$ws.Range("F13").Value = 0.9979078733253837
$ws.Range("G13").Value = 0.9375139589488251
$ws.Range("H13").Value = 1.089588817401622
$ws.Range("I13").Value = 1.019528889404876
$ws.Range("J13").Value = 0.9066257400676264
$ws.Range("K13").Value = 0.9667142791367338
$ws.Range("L13").Value = 0.9823110762310587
$ws.Range("M13").Value = 0.9963280162256956

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8923
$ws.Range("D14").Value = 1.296451999999999
$ws.Range("E14").Value = 0.9335399999999994
$ws.Range("F14").Value = 0.8923
$ws.Range("G14").Value = 1.161671999999999
$ws.Range("H14").Value = 0.8393640000000016
$ws.Range("I14").Value = 0.9183279999999989
$ws.Range("J14").Value = 1.296451999999999
$ws.Range("K14").Value = 1.114995999999999
$ws.Range("L14").Value = 1.003648
$ws.Range("M14").Value = 1.006942666666666

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.8437374999999998
$ws.Range("D15").Value = 1.51
$ws.Range("E15").Value = 0.88
$ws.Range("F15").Value = 0.8437374999999998
$ws.Range("G15").Value = 1.29
$ws.Range("H15").Value = 0.6899999999999999
$ws.Range("I15").Value = 0.86
$ws.Range("J15").Value = 1.51
$ws.Range("K15").Value = 1.195
$ws.Range("L15").Value = 1.01936875
$ws.Range("M15").Value = 1.012289583333333

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9065920592896004
$ws.Range("D16").Value = 1.2962829060096
$ws.Range("E16").Value = 0.9277807521792045
$ws.Range("F16").Value = 0.9065920592896004
$ws.Range("G16").Value = 1.1657767430144
$ws.Range("H16").Value = 0.8170818607104008
$ws.Range("I16").Value = 0.9175305474048041
$ws.Range("J16").Value = 1.2962829060096
$ws.Range("K16").Value = 1.112031829094402
$ws.Range("L16").Value = 1.009311944192001
$ws.Range("M16").Value = 1.005174144768002

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9961450952923423
$ws.Range("D17").Value = 0.9948645301437722
$ws.Range("E17").Value = 0.9953072233463319
$ws.Range("F17").Value = 0.9961450952923423
$ws.Range("G17").Value = 0.9961801600578998
$ws.Range("H17").Value = 0.9944354238425219
$ws.Range("I17").Value = 0.9950813971028302
$ws.Range("J17").Value = 0.9948645301437722
$ws.Range("K17").Value = 0.9950858767450521
$ws.Range("L17").Value = 0.9956154860186972
$ws.Range("M17").Value = 0.9953356382976164

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.001489709722243
$ws.Range("D18").Value = 0.976990827316586
$ws.Range("E18").Value = 1.000190759593898
$ws.Range("F18").Value = 1.001489709722243
$ws.Range("G18").Value = 0.9869174308476603
$ws.Range("H18").Value = 1.005506500119639
$ws.Range("I18").Value = 1.000322646216012
$ws.Range("J18").Value = 0.976990827316586
$ws.Range("K18").Value = 0.988590793455242
$ws.Range("L18").Value = 0.9950402515887424
$ws.Range("M18").Value = 0.995236312302673

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.007815476445126
$ws.Range("D19").Value = 0.9395948733573072
$ws.Range("E19").Value = 1.009061091913582
$ws.Range("F19").Value = 1.007815476445126
$ws.Range("G19").Value = 0.9628406442300881
$ws.Range("H19").Value = 1.034194752031798
$ws.Range("I19").Value = 1.010627183257199
$ws.Range("J19").Value = 0.9395948733573072
$ws.Range("K19").Value = 0.9743279826354446
$ws.Range("L19").Value = 0.9910717295402853
$ws.Range("M19").Value = 0.9940223368725167
